$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestSteps = $wb.Worksheets.Item("Test Steps")

# --- Test Cases sheet -------------------------------------------------
# Runmode/result values: "Yes" -> "yes", "FAIL" -> "PASS"
$wsTestCases.Range("C2").Value = "yes"
$wsTestCases.Range("D3").Value = "PASS"

# --- Test Steps sheet --------------------------------------------------
# The last row that had been run still showed FAIL; fix it to PASS and
# back-fill the still-missing result column for the rest of the steps.
$wsTestSteps.Range("H20").Value = "PASS"
$wsTestSteps.Range("H21").Value = "PASS"
$wsTestSteps.Range("H22").Value = "PASS"
$wsTestSteps.Range("H23").Value = "PASS"
$wsTestSteps.Range("H24").Value = "PASS"
$wsTestSteps.Range("H25").Value = "PASS"
$wsTestSteps.Range("H26").Value = "PASS"
$wsTestSteps.Range("H27").Value = "PASS"
$wsTestSteps.Range("H28").Value = "PASS"

# --- Restore cursor / selection positions ------------------------------
# Move the cell cursor on "Test Cases" without leaving it as the active
# tab (it select()s it, note it back to "Test Steps" afterwards, which
# matches the tabSelected flag staying on "Test Steps").
$wsTestCases.Range("C11").Select() | Out-Null
$wsTestSteps.Range("D14").Select() | Out-Null
